$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 148.5
$ws.Range("I4").Value = 124
$ws.Range("K4").Value = 124
$ws.Range("M4").Value = -10
$ws.Range("H6").Value = 288.83334
$ws.Range("I6").Value = 288.83334
$ws.Range("K6").Value = 866.5000200000001
$ws.Range("M6").Value = -754.5000200000001
$ws.Range("H33").Value = 662.5
$ws.Range("I33").Value = 133.41667
$ws.Range("K33").Value = 133.41667
$ws.Range("M33").Value = 95.58332999999999
$ws.Range("H43").Value = 950
$ws.Range("I43").Value = 900
$ws.Range("J43").Value = 1000
$ws.Range("K43").Value = 900
$ws.Range("L43").Value = 1000
$ws.Range("M43").Value = -831
$ws.Range("N43").Value = -1138
$ws.Range("H86").Value = 6266
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 6266
$ws.Range("K86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("M86").Value = 6266
$ws.Range("N86").Value = -8512
$ws.Range("H89").Value = 6266
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 6266
$ws.Range("K89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("M89").Value = 31330
$ws.Range("N89").Value = -42562
$ws.Range("H113").Value = 19065.666
$ws.Range("I113").Value = 49997
$ws.Range("J113").Value = 3600
$ws.Range("K113").Value = 49997
$ws.Range("L113").Value = 3600
$ws.Range("M113").Value = -46743
$ws.Range("N113").Value = -10108
$ws.Range("H132").Value = 793.2727
$ws.Range("I132").Value = 793.2727
$ws.Range("K132").Value = 2379.8181
$ws.Range("M132").Value = 150.1819
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2980.9
$ws.Range("I2").Value = 2923.2222
$ws.Range("K2").Value = 2923.2222
$ws.Range("M2").Value = -2810.2222
$ws.Range("H32").Value = 3054.15
$ws.Range("I32").Value = 3054.15
$ws.Range("K32").Value = 3054.15
$ws.Range("M32").Value = -2767.15
$ws.Range("H45").Value = 4155.3335
$ws.Range("J45").Value = 4586.4
$ws.Range("L45").Value = 4586.4
$ws.Range("N45").Value = -5340.4
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").ClearContents()
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = 0
$ws.Range("H109").Value = 100000
$ws.Range("J109").Value = 100000
$ws.Range("L109").Value = 100000
$ws.Range("N109").Value = -102774
$ws.Range("H116").Value = 2980.9
$ws.Range("I116").Value = 2923.2222
$ws.Range("K116").Value = 2923.2222
$ws.Range("M116").Value = -629.2222000000002
$ws.Range("H131").Value = 90000
$ws.Range("J131").Value = 90000
$ws.Range("L131").Value = 90000
$ws.Range("N131").Value = -100080
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2980.9
$ws.Range("I3").Value = 2923.2222
$ws.Range("K3").Value = 2923.2222
$ws.Range("M3").Value = -2809.2222
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 4281.4443
$ws.Range("I22").Value = 5014.6924
$ws.Range("J22").Value = 2375
$ws.Range("K22").Value = 5014.6924
$ws.Range("L22").Value = 2375
$ws.Range("M22").Value = -4664.6924
$ws.Range("N22").Value = -3075
$ws.Range("H132").Value = 4559.25
$ws.Range("I132").Value = 4302
$ws.Range("J132").Value = 4988
$ws.Range("K132").Value = 12906
$ws.Range("L132").Value = 14964
$ws.Range("M132").Value = -10376
$ws.Range("N132").Value = -20024
$ws.Range("H134").Value = 1696.1428
$ws.Range("I134").Value = 1612.1666
$ws.Range("K134").Value = 4836.4998
$ws.Range("M134").Value = -2301.4998
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 25018.4
$ws.Range("I11").Value = 36666.668
$ws.Range("J11").Value = 7546
$ws.Range("K11").Value = 110000.004
$ws.Range("L11").Value = 22638
$ws.Range("M11").Value = -109860.004
$ws.Range("N11").Value = -22918
$ws.Range("H68").Value = 3149.842
$ws.Range("I68").Value = 3294.5454
$ws.Range("J68").Value = 2950.875
$ws.Range("K68").Value = 9883.636200000001
$ws.Range("L68").Value = 8852.625
$ws.Range("M68").Value = -9072.636200000001
$ws.Range("N68").Value = -10474.625
$ws.Range("H71").Value = 3149.842
$ws.Range("I71").Value = 3294.5454
$ws.Range("J71").Value = 2950.875
$ws.Range("K71").Value = 29650.9086
$ws.Range("L71").Value = 26557.875
$ws.Range("M71").Value = -25594.9086
$ws.Range("N71").Value = -34669.875
$ws.Range("H75").Value = 3444.2
$ws.Range("J75").Value = 3233.3333
$ws.Range("L75").Value = 9699.999899999999
$ws.Range("N75").Value = -11695.9999
$ws.Range("H78").Value = 3444.2
$ws.Range("J78").Value = 3233.3333
$ws.Range("L78").Value = 29099.9997
$ws.Range("N78").Value = -39083.9997
$ws.Range("H109").Value = 2569.3333
$ws.Range("I109").Value = 83.2
$ws.Range("J109").Value = 15000
$ws.Range("K109").Value = 249.6
$ws.Range("L109").Value = 45000
$ws.Range("M109").Value = 790.4
$ws.Range("N109").Value = -47080
$ws.Range("H114").Value = 307
$ws.Range("I114").Value = 301
$ws.Range("K114").Value = 903
$ws.Range("M114").Value = 2351
$ws.Range("H121").Value = 1170.5
$ws.Range("J121").Value = 2011
$ws.Range("L121").Value = 6033
$ws.Range("N121").Value = -8653
$ws.Range("H131").Value = 849
$ws.Range("J131").Value = 999
$ws.Range("L131").Value = 2997
$ws.Range("N131").Value = -13077
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 204.88889
$ws.Range("I2").Value = 53.125
$ws.Range("J2").Value = 326.3
$ws.Range("K2").Value = 53.125
$ws.Range("L2").Value = 326.3
$ws.Range("M2").Value = 59.875
$ws.Range("N2").Value = -552.3
$ws.Range("H12").Value = 77.5
$ws.Range("I12").Value = 105
$ws.Range("K12").Value = 105
$ws.Range("M12").Value = 35
$ws.Range("H14").Value = 145308.62
$ws.Range("J14").Value = 26928.166
$ws.Range("L14").Value = 26928.166
$ws.Range("N14").Value = -27264.166
$ws.Range("H17").Value = 496.66666
$ws.Range("J17").Value = 612.25
$ws.Range("L17").Value = 612.25
$ws.Range("N17").Value = -948.25
$ws.Range("H19").Value = 1837.5
$ws.Range("J19").Value = 1000
$ws.Range("L19").Value = 1000
$ws.Range("N19").Value = -1576
$ws.Range("H22").Value = 6529.6
$ws.Range("I22").Value = 4466.3335
$ws.Range("J22").Value = 9624.5
$ws.Range("K22").Value = 4466.3335
$ws.Range("L22").Value = 9624.5
$ws.Range("M22").Value = -3937.3335
$ws.Range("N22").Value = -10682.5
$ws.Range("H43").Value = 14249.875
$ws.Range("I43").Value = 29999
$ws.Range("J43").Value = 12000
$ws.Range("K43").Value = 29999
$ws.Range("L43").Value = 12000
$ws.Range("M43").Value = -29848
$ws.Range("N43").Value = -12302
$ws.Range("H107").Value = 11116.5
$ws.Range("I107").Value = 237
$ws.Range("J107").Value = 16556.25
$ws.Range("K107").Value = 237
$ws.Range("L107").Value = 16556.25
$ws.Range("M107").Value = 1683
$ws.Range("N107").Value = -20396.25
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("M119").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 28750
$ws.Range("I2").Value = 5000
$ws.Range("J2").Value = 36666.668
$ws.Range("K2").Value = 5000
$ws.Range("L2").Value = 36666.668
$ws.Range("M2").Value = -4888
$ws.Range("N2").Value = -36890.668
$ws.Range("H7").Value = 5000
$ws.Range("I7").Value = 5000
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 5000
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -4888
$ws.Range("N7").Value = -5224
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").ClearContents()
$ws.Range("N18").Value = 0
$ws.Range("H22").Value = 2201.4211
$ws.Range("I22").Value = 1610.6666
$ws.Range("K22").Value = 1610.6666
$ws.Range("M22").Value = -1315.6666
$ws.Range("H27").Value = 2201.4211
$ws.Range("I27").Value = 1610.6666
$ws.Range("K27").Value = 1610.6666
$ws.Range("M27").Value = -1503.6666
$ws.Range("H46").Value = 4041.6667
$ws.Range("I46").Value = 1750
$ws.Range("K46").Value = 1750
$ws.Range("M46").Value = -1562
$ws.Range("H55").Value = 565.1111
$ws.Range("I55").Value = 428.375
$ws.Range("J55").Value = 674.5
$ws.Range("K55").Value = 428.375
$ws.Range("L55").Value = 674.5
$ws.Range("M55").Value = -255.375
$ws.Range("N55").Value = -1020.5
$ws.Range("H68").Value = 5251.5
$ws.Range("I68").Value = 500
$ws.Range("K68").Value = 500
$ws.Range("M68").Value = 249
$ws.Range("H71").Value = 5251.5
$ws.Range("I71").Value = 500
$ws.Range("K71").Value = 2500
$ws.Range("M71").Value = 1244
$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 5000
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 15000
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -12530
$ws.Range("N126").Value = -19940
$ws.Range("H132").Value = 2000
$ws.Range("J132").Value = 2000
$ws.Range("L132").Value = 6000
$ws.Range("N132").Value = -11060
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2501.5
$ws.Range("J62").Value = 2501.5
$ws.Range("L62").Value = 2501.5
$ws.Range("N62").Value = -3749.5
$ws.Range("H65").Value = 2501.5
$ws.Range("J65").Value = 2501.5
$ws.Range("L65").Value = 12507.5
$ws.Range("N65").Value = -18747.5
$ws.Range("H106").Value = 40000
$ws.Range("J106").Value = 40000
$ws.Range("L106").Value = 40000
$ws.Range("N106").Value = -42524

Write-Output "Applied all changes"